# "first year model artifact removal"
# The mass-balance model was re-run after removing the first year of data
# from the underlying model artifact; this changed the hard-coded
# mass-balance figures on the "5) Mass balances" sheet (B3:F7), which in
# turn ripple through the dependent percentage formulas (B10:F14) and the
# SUM() totals in column G automatically on recalculation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 2) Model parameters: re-enter the "1 - x" formulas as a single fill
# across C27:G28 so Excel stores them as one shared-formula group
# (matches a select + fill-down/right re-entry of the formula).
# ---------------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("2) Model parameters")
$wsParams.Range("C27:G28").Formula = "=1-C23"

# ---------------------------------------------------------------------
# 5) Mass balances: updated model output values (B3:F7)
# ---------------------------------------------------------------------
$wsMB = $wb.Worksheets.Item("5) Mass balances")

$wsMB.Range("B3").Value = 39.911058312097197
$wsMB.Range("C3").Value = 32.9424128896
$wsMB.Range("D3").Value = -48.759044861748599
$wsMB.Range("E3").Value = -5.3017393251846103
$wsMB.Range("F3").Value = -17.415200198005

$wsMB.Range("B4").Value = 66.303159712552102
$wsMB.Range("C4").Value = 54.7168588588939
$wsMB.Range("D4").Value = -16.5182492312344
$wsMB.Range("E4").Value = -44.3848619195685
$wsMB.Range("F4").Value = -61.055035308213

$wsMB.Range("B5").Value = 72.775175448556794
$wsMB.Range("C5").Value = 11.4771511916558
$wsMB.Range("D5").Value = -24.0331239756059
$wsMB.Range("E5").Value = -6.3727021573982903
$wsMB.Range("F5").Value = -53.415329388723002

$wsMB.Range("B6").Value = 13.769777258250601
$wsMB.Range("C6").Value = 28.140834484649201
$wsMB.Range("D6").Value = -37.084408704860003
$wsMB.Range("E6").Value = -1.21205080195923
$wsMB.Range("F6").Value = -4.1299826785709701

$wsMB.Range("B7").Value = 39.576465329266298
$wsMB.Range("C7").Value = 26.4145231974755
$wsMB.Range("D7").Value = -33.375631727461503
$wsMB.Range("E7").Value = -20.100651477426499
$wsMB.Range("F7").Value = -12.234891490866801

# Row 7 (the "5) Total" row) loses its distinct bordered style and
# switches to the same plain numeric style as the rows above it once the
# new values are entered.
$wsMB.Range("B7:F7").Borders.LineStyle = -4142

# ---------------------------------------------------------------------
# Tab / selection: the workbook now opens on "5) Mass balances" (was
# "3)Equations"), with the cursor left on D18.
# ---------------------------------------------------------------------
$wsMB.Activate()
$wsMB.Range("D18").Select()
